# Insert a new row for "Algeria" (cc=612) right after the header/Albania rows,
# i.e. as the new row 3, pushing Angola and everything below it down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = 612
$ws.Range("B3").Value = "Algeria"

# Update selection to reflect the author's new cursor position.
$ws.Range("A4").Select()
